$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 5615.6665
$ws.Range("I98").Value = 5738.8667
$ws.Range("K98").Value = 5738.8667
$ws.Range("M98").Value = -4240.8667
$ws.Range("H100").Value = 2044.5
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H101").Value = 466.66666
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("H122").Value = 5615.6665
$ws.Range("I122").Value = 5738.8667
$ws.Range("K122").Value = 17216.6001
$ws.Range("M122").Value = -14766.6001
$ws.Range("H131").Value = 3336061.2
$ws.Range("I131").Value = 3336061.2
$ws.Range("K131").Value = 10008183.6
$ws.Range("M131").Value = -10003143.6
$ws.Range("H137").Value = 3998.4878
$ws.Range("I137").Value = 1512.24
$ws.Range("J137").Value = 7883.25
$ws.Range("K137").Value = 4536.72
$ws.Range("L137").Value = 23649.75
$ws.Range("M137").Value = -1986.72
$ws.Range("N137").Value = -28749.75
$ws.Range("H138").Value = 256280.5
$ws.Range("I138").Value = 3668.5
$ws.Range("J138").Value = 400630.22
$ws.Range("K138").Value = 11005.5
$ws.Range("L138").Value = 1201890.66
$ws.Range("M138").Value = -5865.5
$ws.Range("N138").Value = -1212170.66

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 695.5
$ws.Range("I2").Value = 718.08
$ws.Range("J2").Value = 614.8570999999999
$ws.Range("K2").Value = 718.08
$ws.Range("L2").Value = 614.8570999999999
$ws.Range("M2").Value = -605.08
$ws.Range("N2").Value = -840.8570999999999
$ws.Range("H32").Value = 4016.1567
$ws.Range("I32").Value = 2894.052
$ws.Range("K32").Value = 2894.052
$ws.Range("M32").Value = -2607.052
$ws.Range("H63").Value = 2436.25
$ws.Range("I63").Value = 2436.25
$ws.Range("K63").Value = 2436.25
$ws.Range("M63").Value = -1750.25
$ws.Range("H66").Value = 2436.25
$ws.Range("I66").Value = 2436.25
$ws.Range("K66").Value = 12181.25
$ws.Range("M66").Value = -8749.25
$ws.Range("H116").Value = 695.5
$ws.Range("I116").Value = 718.08
$ws.Range("J116").Value = 614.8570999999999
$ws.Range("K116").Value = 718.08
$ws.Range("L116").Value = 614.8570999999999
$ws.Range("M116").Value = 1575.92
$ws.Range("N116").Value = -5202.8571
$ws.Range("H121").Value = 49252.5
$ws.Range("J121").Value = 49252.5
$ws.Range("L121").Value = 49252.5
$ws.Range("N121").Value = -52746.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 695.5
$ws.Range("I3").Value = 718.08
$ws.Range("J3").Value = 614.8570999999999
$ws.Range("K3").Value = 718.08
$ws.Range("L3").Value = 614.8570999999999
$ws.Range("M3").Value = -604.08
$ws.Range("N3").Value = -842.8570999999999
$ws.Range("H20").Value = 26885196
$ws.Range("I20").Value = 33336750
$ws.Range("K20").Value = 33336750
$ws.Range("M20").Value = -33336503
$ws.Range("H29").Value = 8700
$ws.Range("I29").Value = 8700
$ws.Range("K29").Value = 8700
$ws.Range("M29").Value = -8411

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4234.909
$ws.Range("I86").Value = 4088.5
$ws.Range("K86").Value = 4088.5
$ws.Range("M86").Value = -2965.5
$ws.Range("H89").Value = 4234.909
$ws.Range("I89").Value = 4088.5
$ws.Range("K89").Value = 20442.5
$ws.Range("M89").Value = -14826.5
$ws.Range("H99").Value = 4613.6113
$ws.Range("I99").Value = 4758.077
$ws.Range("J99").Value = 4238
$ws.Range("K99").Value = 4758.077
$ws.Range("L99").Value = 4238
$ws.Range("M99").Value = -3260.077
$ws.Range("N99").Value = -7234
$ws.Range("H126").Value = 4613.6113
$ws.Range("I126").Value = 4758.077
$ws.Range("J126").Value = 4238
$ws.Range("K126").Value = 14274.231
$ws.Range("L126").Value = 12714
$ws.Range("M126").Value = -11804.231
$ws.Range("N126").Value = -17654
$ws.Range("H132").Value = 2305.082
$ws.Range("I132").Value = 1751.56
$ws.Range("J132").Value = 4821.091
$ws.Range("K132").Value = 5254.68
$ws.Range("L132").Value = 14463.273
$ws.Range("M132").Value = -2724.68
$ws.Range("N132").Value = -19523.273
$ws.Range("H134").Value = 2355.2334
$ws.Range("I134").Value = 2230.5715
$ws.Range("K134").Value = 6691.7145
$ws.Range("M134").Value = -4156.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1250445.5
$ws.Range("J97").Value = 594
$ws.Range("L97").Value = 1782
$ws.Range("N97").Value = -2774
$ws.Range("H137").Value = 4002.3044
$ws.Range("J137").Value = 3891
$ws.Range("L137").Value = 11673
$ws.Range("N137").Value = -21873
$ws.Range("H140").Value = 11962.462
$ws.Range("J140").Value = 21099.9
$ws.Range("L140").Value = 63299.7
$ws.Range("N140").Value = -73659.70000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 7777
$ws.Range("J20").Value = 7777
$ws.Range("L20").Value = 7777
$ws.Range("N20").Value = -8267
$ws.Range("H24").Value = 10916.375
$ws.Range("I24").Value = 10996.5
$ws.Range("J24").Value = 10909.091
$ws.Range("K24").Value = 10996.5
$ws.Range("L24").Value = 10909.091
$ws.Range("M24").Value = -10823.5
$ws.Range("N24").Value = -11255.091
$ws.Range("H70").Value = 14772489
$ws.Range("J70").Value = 14363
$ws.Range("L70").Value = 14363
$ws.Range("N70").Value = -14903
$ws.Range("H73").Value = 14772489
$ws.Range("J73").Value = 14363
$ws.Range("L73").Value = 14363
$ws.Range("N73").Value = -16235

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6023.75
$ws.Range("I7").Value = 4469.143
$ws.Range("J7").Value = 8200.200000000001
$ws.Range("K7").Value = 4469.143
$ws.Range("L7").Value = 8200.200000000001
$ws.Range("M7").Value = -4357.143
$ws.Range("N7").Value = -8424.200000000001
$ws.Range("H22").Value = 514.4375
$ws.Range("I22").Value = 471.2
$ws.Range("J22").Value = 586.5
$ws.Range("K22").Value = 471.2
$ws.Range("L22").Value = 586.5
$ws.Range("M22").Value = -176.2
$ws.Range("N22").Value = -1176.5
$ws.Range("H27").Value = 514.4375
$ws.Range("I27").Value = 471.2
$ws.Range("J27").Value = 586.5
$ws.Range("K27").Value = 471.2
$ws.Range("L27").Value = 586.5
$ws.Range("M27").Value = -364.2
$ws.Range("N27").Value = -800.5
$ws.Range("H40").Value = 5531.1562
$ws.Range("I40").Value = 5546.5356
$ws.Range("J40").Value = 5423.5
$ws.Range("K40").Value = 5546.5356
$ws.Range("L40").Value = 5423.5
$ws.Range("M40").Value = -5410.5356
$ws.Range("N40").Value = -5695.5
$ws.Range("H126").Value = 6023.75
$ws.Range("I126").Value = 4469.143
$ws.Range("J126").Value = 8200.200000000001
$ws.Range("K126").Value = 13407.429
$ws.Range("L126").Value = 24600.6
$ws.Range("M126").Value = -10937.429
$ws.Range("N126").Value = -29540.6
$ws.Range("H132").Value = 5724.143
$ws.Range("I132").Value = 2942.3572
$ws.Range("K132").Value = 8827.071599999999
$ws.Range("M132").Value = -6297.071599999999
$ws.Range("H136").Value = 4142.9375
$ws.Range("I136").Value = 3391.524
$ws.Range("J136").Value = 5577.4546
$ws.Range("K136").Value = 10174.572
$ws.Range("L136").Value = 16732.3638
$ws.Range("M136").Value = -7624.572
$ws.Range("N136").Value = -21832.3638
$ws.Range("H139").Value = 69998.91
$ws.Range("I139").Value = 69998
$ws.Range("K139").Value = 69998
$ws.Range("M139").Value = -64858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 25000
$ws.Range("J25").Value = 25000
$ws.Range("L25").Value = 25000
$ws.Range("N25").Value = -25586
$ws.Range("H126").Value = 2554.8333
$ws.Range("I126").Value = 1380.4286
$ws.Range("J126").Value = 4199
$ws.Range("K126").Value = 4141.2858
$ws.Range("L126").Value = 12597
$ws.Range("M126").Value = -1671.2858
$ws.Range("N126").Value = -17537
$ws.Range("H137").Value = 78180.91
$ws.Range("J137").Value = 78180.91
$ws.Range("L137").Value = 78180.91
$ws.Range("N137").Value = -88380.91
